$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two data rows (old rows 6 and 7), shifting remaining rows up.
$ws.Range("6:7").Delete() | Out-Null

# New data for rows 2-5 (Sending cluster, Ligand symbol, Receptor symbol, Target cluster, ...)
$data = @(
    @("FAPs","Efna5","Ephb1","ECs",3,1,2.900731333333333,8.702194,0.8130494232775288,0.8130494232775289,3,1,1.475208,4.425624,0.8210007041987012,0.8210007041987013,4.279182068784,38.512638619056,0.667514149059199,0.6675141490591993),
    @("FAPs","Efna5","Ephb1","MuSCs",3,1,2.900731333333333,8.702194,0.8130494232775288,0.8130494232775289,3,1,0.3216333333333333,0.9649,0.1789992958012987,0.1789992958012987,0.9329718878444444,8.396746990600001,0.1455352742183297,0.1455352742183297),
    @("MuSCs","Efna5","Ephb1","ECs",3,1,0.6669870000000001,2.000961,0.1869505767224711,0.1869505767224711,3,1,1.475208,4.425624,0.8210007041987012,0.8210007041987013,0.9839445582960002,8.855501024664001,0.1534865551395021,0.1534865551395021),
    @("MuSCs","Efna5","Ephb1","MuSCs",3,1,0.6669870000000001,2.000961,0.1869505767224711,0.1869505767224711,3,1,0.3216333333333333,0.9649,0.1789992958012987,0.1789992958012987,0.2145252521,1.9307272689,0.03346402158296899,0.033464021582969)
)

$rowIdx = 2
foreach ($row in $data) {
    $colIdx = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIdx, $colIdx).Value = $val
        $colIdx++
    }
    $rowIdx++
}
